# Generate Report for Archive
#
# Update the localization status from "Ready for handoff" to "In Translation"
# across every sheet that carries the status column, then let the column
# widths re-flow to fit the (shorter) new text, matching the narrower
# "Status"/"zh-cn"/"de-de" columns produced when the report was regenerated.

$wb = $excel.ActiveWorkbook

# 1) Update the status text everywhere it appears (Overview sheet's per-language
#    status columns, plus each language sheet's "Status" column).
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation") | Out-Null
}

# 2) Re-fit the columns that held the status text so their width reflects the
#    shorter replacement string.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5   # "zh-cn" status column
$overview.Columns.Item(6).ColumnWidth = 12.5   # "de-de" status column

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # "Status" column

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5        # "Status" column
